# Appends the new daily rolling-7-day rows (comuni MO) for
# 2021-09-02 .. 2021-09-09 (serials 44441-44448), matching the
# existing "date in col A (style copied from the row above) + 49
# numeric comuni columns B:AX" layout used by every prior row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per array, values in column order A, B, C ... AX.
$newRows = @(
    @(44441,14,10,5,118,43,6,14,6,5,5,24,27,10,1,1,0,15,4,1,13,176,2,21,6,21,3,0,11,14,12,13,35,12,13,3,5,12,39,1,34,786,2,7,0,0,2,10,5,5),
    @(44442,12,10,5,112,33,8,11,7,6,5,30,29,8,1,2,0,14,4,1,19,187,4,19,5,26,4,0,10,12,7,8,37,14,13,4,4,9,37,1,32,775,2,3,0,0,1,11,4,4),
    @(44443,16,7,5,105,35,8,10,7,6,4,32,33,9,1,2,0,14,3,0,22,183,6,17,7,24,5,1,8,10,6,8,36,15,9,3,6,9,38,1,27,750,2,1,0,0,1,6,0,2),
    @(44444,9,7,0,103,28,8,9,6,6,3,25,26,7,1,2,0,16,3,0,18,186,6,14,6,19,5,1,7,9,7,8,36,14,4,3,7,9,30,1,25,683,2,0,0,0,0,6,0,1),
    @(44445,9,5,0,76,28,9,7,7,6,4,25,25,8,1,2,0,13,2,0,23,215,5,9,10,19,4,1,4,7,7,6,32,15,1,4,8,13,26,0,20,666,3,0,0,0,0,6,0,1),
    @(44446,8,7,0,72,31,7,6,7,4,1,18,15,8,1,2,0,13,2,0,24,205,5,9,10,17,3,1,6,6,5,4,35,16,2,4,9,18,23,0,15,630,3,0,0,0,0,5,0,3),
    @(44447,7,7,0,74,32,7,9,6,4,1,18,15,8,1,2,0,13,2,0,25,200,6,8,10,17,3,1,6,6,6,4,34,15,2,4,9,18,23,0,16,630,3,0,0,0,0,5,0,3),
    @(44448,6,6,1,70,27,7,8,7,1,1,19,13,6,0,2,0,9,1,2,26,196,6,2,10,31,4,1,7,6,7,5,29,15,4,2,7,16,26,0,13,608,3,0,0,0,0,1,0,5),
)

$startRow = 367
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $vals = $newRows[$i]
    $arr = New-Object "object[,]" 1,$vals.Count
    for ($c = 0; $c -lt $vals.Count; $c++) { $arr[0,$c] = $vals[$c] }
    $r = $startRow + $i
    $ws.Range("A$r" + ":AX$r").Value = $arr
}

# Column A carries the same centered/bordered date style ("s=2"
# in the sheet XML) as every other date cell; copy it down from
# the last pre-existing row (366) instead of re-deriving the xf.
$ws.Range("A366").Copy() | Out-Null
$ws.Range("A367:A374").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
